# Regenerate the "K" column (column G) values for save_data/rogers_taylor.
# The workbook's column headers are: B=date, C=TB, D=PC, E=dS0, F=dSF,
# G=K, H=IP, I=I0, J=IF. This script rewrites the recomputed K values
# (formerly a "Strike#" style count) for each data row, row 2 through 70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 1
    17 = 3
    18 = 2
    19 = 0
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 2
    31 = 1
    32 = 0
    33 = 1
    34 = 2
    35 = 2
    36 = 2
    37 = 2
    38 = 3
    39 = 2
    40 = 3
    41 = 1
    42 = 3
    43 = 2
    44 = 1
    45 = 2
    46 = 1
    47 = 0
    48 = 1
    49 = 0
    50 = 2
    51 = 1
    52 = 1
    54 = 1
    55 = 1
    56 = 0
    57 = 0
    58 = 0
    59 = 1
    60 = 0
    61 = 0
    62 = 1
    64 = 3
    65 = 1
    66 = 1
    67 = 2
    68 = 1
    69 = 2
    70 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
